$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "34.110.67"
$ws.Cells.Item(2, 5).Value = "  -0.69%  "
$ws.Cells.Item(3, 4).Value = "1.784.20"
$ws.Cells.Item(3, 5).Value = "  -2.89%  "
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.00"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(5, 4).Value = "224.72"
$ws.Cells.Item(5, 5).Value = "  -0.45%  "
$ws.Cells.Item(6, 5).Value = "  -1.31%  "
$ws.Cells.Item(7, 5).Value = "  +0.12%  "
$ws.Cells.Item(8, 5).Value = "  +2.11%  "
$ws.Cells.Item(9, 5).Value = "  -2.46%  "
$ws.Cells.Item(10, 4).Value = "0.0711"
$ws.Cells.Item(10, 5).Value = "  -1.84%  "
$ws.Cells.Item(11, 4).Value = "0.0936"
$ws.Cells.Item(11, 5).Value = "  +0.53%  "
$ws.Cells.Item(12, 4).Value = "2.042.43"
$ws.Cells.Item(12, 5).Value = "  -2.84%  "
$ws.Cells.Item(13, 4).Value = "10.99"
$ws.Cells.Item(13, 5).Value = "  +1.80%  "
$ws.Cells.Item(14, 4).Value = "1.789.73"
$ws.Cells.Item(14, 5).Value = "  -2.48%  "
$ws.Cells.Item(15, 2).Value = "WrappedBTC"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Cells.Item(15, 4).Value = "34.051.61"
$ws.Cells.Item(15, 5).Value = "  -0.94%  "
$ws.Cells.Item(16, 2).Value = "Polygon"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Cells.Item(16, 4).Value = "0.621"
$ws.Cells.Item(16, 5).Value = "  -4.03%  "
$ws.Cells.Item(17, 5).Value = "  -4.42%  "
$ws.Cells.Item(18, 4).Value = "67.78"
$ws.Cells.Item(18, 5).Value = "  -3.05%  "
$ws.Cells.Item(19, 4).Value = "245.05"
$ws.Cells.Item(19, 5).Value = "  -2.92%  "
$ws.Cells.Item(20, 4).Value = "0.0₃0789"
$ws.Cells.Item(20, 5).Value = "  -1.18%  "
$ws.Cells.Item(21, 5).Value = "  +0.14%  "
$ws.Cells.Item(22, 4).Value = "10.82"
$ws.Cells.Item(22, 5).Value = "  -3.73%  "
$ws.Cells.Item(23, 4).Value = "4.08"
$ws.Cells.Item(23, 5).Value = "  -4.73%  "
$ws.Cells.Item(24, 5).Value = "  -3.03%  "
$ws.Cells.Item(25, 4).Value = "160.74"
$ws.Cells.Item(26, 4).Value = "16.32"
$ws.Cells.Item(26, 5).Value = "  -2.55%  "
$ws.Cells.Item(27, 4).Value = "7.06"
$ws.Cells.Item(27, 5).Value = "  -3.01%  "
$ws.Cells.Item(28, 4).Value = "0.112"
$ws.Cells.Item(28, 5).Value = "  -2.67%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "1.00"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +0.15%  "
$ws.Cells.Item(30, 2).Value = "PancakeSwap"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(30, 4).Value = "1.21"
$ws.Cells.Item(30, 5).Value = "  -0.27%  "
$ws.Cells.Item(31, 2).Value = "Hedera"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(31, 4).Value = "0.0512"
$ws.Cells.Item(31, 5).Value = "  -4.78%  "
$ws.Cells.Item(32, 4).Value = "3.64"
$ws.Cells.Item(32, 5).Value = "  -4.52%  "
$ws.Cells.Item(33, 5).Value = "  -2.33%  "
$ws.Cells.Item(34, 5).Value = "  -5.41%  "
$ws.Cells.Item(35, 4).Value = "1.393.60"
$ws.Cells.Item(35, 5).Value = "  -3.96%  "
$ws.Cells.Item(36, 5).Value = "  -0.97%  "
$ws.Cells.Item(37, 5).Value = "  -1.92%  "
$ws.Cells.Item(38, 5).Value = "  -3.04%  "
$ws.Cells.Item(39, 2).Value = "HuobiToken"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(39, 4).Value = "2.35"
$ws.Cells.Item(39, 5).Value = "  -0.16%  "
$ws.Cells.Item(40, 2).Value = "RenderToken"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "2.20"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +2.67%  "
$ws.Cells.Item(41, 5).Value = "  -5.80%  "
$ws.Cells.Item(42, 5).Value = "  -2.83%  "
$ws.Cells.Item(43, 4).Value = "78.13"
$ws.Cells.Item(43, 5).Value = "  -4.88%  "
$ws.Cells.Item(44, 4).Value = "0.0₆0140"
$ws.Cells.Item(44, 5).Value = "  +12.08%  "
$ws.Cells.Item(45, 5).Value = "  +3.02%  "
$ws.Cells.Item(46, 4).Value = "108.33"
$ws.Cells.Item(46, 5).Value = "  +1.35%  "
$ws.Cells.Item(47, 5).Value = "  -0.80%  "
$ws.Cells.Item(48, 4).Value = "12.39"
$ws.Cells.Item(48, 5).Value = "  +3.91%  "
$ws.Cells.Item(49, 4).Value = "5.84"
$ws.Cells.Item(49, 5).Value = "  -4.28%  "
$ws.Cells.Item(50, 4).Value = "1.941.95"
$ws.Cells.Item(50, 5).Value = "  -2.69%  "
$ws.Cells.Item(51, 5).Value = "  +0.08%  "
